$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.2489575877753215
$ws.Range("D2").Value = 0.8057030807601888

# Row 3
$ws.Range("C3").Value = -0.3483473062333099
$ws.Range("D3").Value = 0.730890740844286

# Row 4
$ws.Range("C4").Value = -2.134569066848549
$ws.Range("D4").Value = 0.04418538521563287

# Row 5
$ws.Range("C5").Value = -1.829870204051465
$ws.Range("D5").Value = 0.08085396584447979
$ws.Range("G5").Value = "No"

# Row 6
$ws.Range("C6").Value = -0.2389559571122948
$ws.Range("D6").Value = 0.8133516985237434

# Row 7
$ws.Range("C7").Value = -1.737557501603092
$ws.Range("D7").Value = 0.09626981989854833

# Row 8
$ws.Range("C8").Value = -2.006773555597452
$ws.Range("D8").Value = 0.05721786415425756
$ws.Range("G8").Value = "No"

# Row 9
$ws.Range("C9").Value = -1.581004099212165
$ws.Range("D9").Value = 0.1281484256702587

# Row 10
$ws.Range("C10").Value = -1.865178419927265
$ws.Range("D10").Value = 0.07555025246441782
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = -0.3396166242535306
$ws.Range("D11").Value = 0.7373661014924631
